$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the width currently used by column D (shared across D:K) so the
# two freshly inserted columns can be sized to match once they exist.
$dataColWidth = $ws.Range("D1").EntireColumn.ColumnWidth

# Insert two new columns before column D (D:E). This shifts the existing
# D:K quarterly data two columns to the right (becomes F:M).
$ws.Range("D1:E1").EntireColumn.Insert()

# The newly inserted D:E columns start out with default/blank formatting.
# Copy number formats/styles from F:G (which now hold the original D:E
# data) so the new columns look like the rest of the quarter columns -
# date format on the "Period Ending" rows (7/38/80), #,##0 elsewhere.
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the column width of the other quarter columns.
$ws.Range("D1:E1").EntireColumn.ColumnWidth = $dataColWidth

# Populate the new D:E columns with the newest two reported quarters.
$ws.Range("D7").Value = 43498
$ws.Range("E7").Value = 43407
$ws.Range("D8").Value = 14801000
$ws.Range("E8").Value = 9590000
$ws.Range("D9").Value = 11518000
$ws.Range("E9").Value = 7266000
$ws.Range("D10").Value = 3283000
$ws.Range("E10").Value = 2324000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = -1000
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 13823000
$ws.Range("E17").Value = 9268000
$ws.Range("D18").Value = 978000
$ws.Range("E18").Value = 322000
$ws.Range("D20").Value = 14000
$ws.Range("E20").Value = 23000
$ws.Range("D21").Value = 1212000
$ws.Range("E21").Value = 537000
$ws.Range("D22").Value = 20000
$ws.Range("E22").Value = 15000
$ws.Range("D23").Value = 972000
$ws.Range("E23").Value = 330000
$ws.Range("D24").Value = 237000
$ws.Range("E24").Value = 53000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 735000
$ws.Range("E26").Value = 277000
$ws.Range("D27").Value = 735000
$ws.Range("E27").Value = 277000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -14000
$ws.Range("E32").Value = -23000
$ws.Range("D33").Value = 735000
$ws.Range("E33").Value = 277000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 735000
$ws.Range("E35").Value = 277000
$ws.Range("D38").Value = 43498
$ws.Range("E38").Value = 43407
$ws.Range("D41").Value = 1980000
$ws.Range("E41").Value = 1228000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 76000
$ws.Range("D43").Value = 1015000
$ws.Range("E43").Value = 921000
$ws.Range("D44").Value = 5409000
$ws.Range("E44").Value = 8168000
$ws.Range("D45").Value = 466000
$ws.Range("E45").Value = 508000
$ws.Range("D46").Value = 8870000
$ws.Range("E46").Value = 10901000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2510000
$ws.Range("E48").Value = 2525000
$ws.Range("D49").Value = 933000
$ws.Range("E49").Value = 939000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 588000
$ws.Range("E52").Value = 635000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 12901000
$ws.Range("E54").Value = 15000000
$ws.Range("D57").Value = 5257000
$ws.Range("E57").Value = 7964000
$ws.Range("D58").Value = 56000
$ws.Range("E58").Value = 46000
$ws.Range("D59").Value = 2200000
$ws.Range("E59").Value = 1923000
$ws.Range("D60").Value = 7513000
$ws.Range("E60").Value = 9933000
$ws.Range("D61").Value = 1332000
$ws.Range("E61").Value = 1280000
$ws.Range("D62").Value = 750000
$ws.Range("E62").Value = 775000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 9595000
$ws.Range("E66").Value = 11988000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 2985000
$ws.Range("E72").Value = 2685000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3306000
$ws.Range("E76").Value = 3012000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43498
$ws.Range("E80").Value = 43407
$ws.Range("D81").Value = 735000
$ws.Range("E81").Value = 277000
$ws.Range("D83").Value = 220000
$ws.Range("E83").Value = 192000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1301000
$ws.Range("E89").Value = -1000
$ws.Range("D91").Value = -200000
$ws.Range("E91").Value = -244000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -66000
$ws.Range("E94").Value = -626000
$ws.Range("D96").Value = -121000
$ws.Range("E96").Value = -123000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -492000
$ws.Range("E100").Value = -2000
$ws.Range("D101").Value = 2000
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 745000
$ws.Range("E102").Value = -629000

Write-Host "Inserted 2 columns and populated new quarterly data."
